$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 "51.292.59"
Set-TextValue 2 5 "  -0.11%  "

Set-TextValue 3 4 "2.971.51"
Set-TextValue 3 5 "  +1.92%  "

Set-TextValue 4 4 "0.999"
Set-TextValue 4 5 "  -0.04%  "

Set-TextValue 5 4 "382.38"
Set-TextValue 5 5 "  +1.78%  "

Set-TextValue 6 4 "102.75"
Set-TextValue 6 5 "  +0.18%  "

Set-TextValue 7 5 "  -0.39%  "

Set-TextValue 8 5 "  -0.01%  "

Set-TextValue 9 5 "  +1.33%  "

Set-TextValue 10 4 "36.61"
Set-TextValue 10 5 "  -0.77%  "

Set-TextValue 11 5 "  +0.04%  "

Set-TextValue 12 4 "0.0841"
Set-TextValue 12 5 "  +1.00%  "

Set-TextValue 13 4 "3.445.91"
Set-TextValue 13 5 "  +1.84%  "

Set-TextValue 14 4 "18.13"
Set-TextValue 14 5 "  -0.73%  "

Set-TextValue 15 5 "  +2.03%  "

Set-TextValue 16 4 "2.972.69"
Set-TextValue 16 5 "  +1.79%  "

Set-TextValue 17 4 "0.990"
Set-TextValue 17 5 "  +7.08%  "

Set-TextValue 18 4 "51.228.35"
Set-TextValue 18 5 "  -0.10%  "

Set-TextValue 19 5 "  -4.31%  "

Set-TextValue 20 5 "  +0.21%  "

Set-TextValue 21 5 "  -1.54%  "

Set-TextValue 22 4 "0.0₃0956"
Set-TextValue 22 5 "  +1.43%  "

Set-TextValue 23 4 "68.78"
Set-TextValue 23 5 "  +0.78%  "

Set-TextValue 24 4 "262.07"
Set-TextValue 24 5 "  +0.30%  "

Set-TextValue 25 5 "  +5.28%  "

Set-TextValue 26 4 "8.13"
Set-TextValue 26 5 "  +12.91%  "

Set-TextValue 27 4 "7.56"
Set-TextValue 27 5 "  +11.23%  "

Set-TextValue 28 5 "  +12.75%  "

Set-TextValue 29 2 "LEO"
Set-TextValue 29 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue 29 4 "4.10"
Set-TextValue 29 5 "  -0.45%  "

Set-TextValue 30 2 "Kaspa"
Set-TextValue 30 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue 30 4 "0.167"
Set-TextValue 30 5 "  -1.44%  "

Set-TextValue 31 4 "0.999"
Set-TextValue 31 5 "  -0.10%  "

Set-TextValue 32 5 "  +0.70%  "

Set-TextValue 33 4 "9.84"
Set-TextValue 33 5 "  +0.77%  "

Set-TextValue 34 2 "OKB"
Set-TextValue 34 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue 34 4 "50.92"
Set-TextValue 34 5 "  -0.66%  "

Set-TextValue 35 2 "InjectiveProtocol"
Set-TextValue 35 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 35 4 "34.10"
Set-TextValue 35 5 "  +0.62%  "

Set-TextValue 36 5 "  -2.59%  "

Set-TextValue 37 5 "  +6.17%  "

Set-TextValue 38 5 "  -0.13%  "

Set-TextValue 39 4 "2.99"
Set-TextValue 39 5 "  +0.21%  "

Set-TextValue 40 4 "17.05"
Set-TextValue 40 5 "  +1.09%  "

Set-TextValue 41 4 "2.55"
Set-TextValue 41 5 "  +0.81%  "

Set-TextValue 42 5 "  +1.61%  "

Set-TextValue 43 5 "  -0.81%  "

Set-TextValue 44 4 "122.55"
Set-TextValue 44 5 "  -0.04%  "

Set-TextValue 45 4 "21.27"
Set-TextValue 45 5 "  -0.95%  "

Set-TextValue 46 5 "  +0.09%  "

Set-TextValue 47 5 "  +2.34%  "

Set-TextValue 49 4 "2.023.93"
Set-TextValue 49 5 "  +0.10%  "

Set-TextValue 50 4 "3.24"
Set-TextValue 50 5 "  +3.23%  "

Set-TextValue 51 4 "0.0335"
Set-TextValue 51 5 "  +6.43%  "
